$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column G (Quantidade),
# shifting Quantidade/Data Venda/Hora Venda right from G:I to I:K.
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

# Keep ID Venda (col A), Data Venda (col J) and Hora Venda (col K) as
# plain text so Excel does not reinterpret numeric-looking / date-looking
# strings as numbers or dates.
$ws.Columns.Item(1).NumberFormat = "@"
$ws.Columns.Item(10).NumberFormat = "@"
$ws.Columns.Item(11).NumberFormat = "@"

# New header cells for the inserted columns.
$ws.Cells.Item(1, 7).Value = "Desconto"
$ws.Cells.Item(1, 8).Value = "Desc. Uni."

# Row 2 -> ID 49 (originally row 6's data)
$ws.Cells.Item(2, 1).Value = "49"
$ws.Cells.Item(2, 6).Value = 425.7
$ws.Cells.Item(2, 7).Value = 100
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 100
$ws.Cells.Item(2, 10).Value = "08/02/2023"
$ws.Cells.Item(2, 11).Value = "10:30:42.000"

# Row 3 -> ID 48 (originally row 5's data)
$ws.Cells.Item(3, 1).Value = "48"
$ws.Cells.Item(3, 6).Value = 425.7
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(3, 10).Value = "08/02/2023"
$ws.Cells.Item(3, 11).Value = "09:45:15.000"

# Row 4 -> ID 54 (stays in place, gains Desconto / Desc. Uni. values)
$ws.Cells.Item(4, 1).Value = "54"
$ws.Cells.Item(4, 6).Value = 425.7
$ws.Cells.Item(4, 7).Value = 100
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 100
$ws.Cells.Item(4, 10).Value = "08/02/2023"
$ws.Cells.Item(4, 11).Value = "11:00:39.000"

# Row 5 -> ID 56 (originally row 3's data)
$ws.Cells.Item(5, 1).Value = "56"
$ws.Cells.Item(5, 6).Value = 310.76
$ws.Cells.Item(5, 7).Value = 73
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 73
$ws.Cells.Item(5, 10).Value = "10/02/2023"
$ws.Cells.Item(5, 11).Value = "13:38:40.000"

# Row 6 -> ID 55 (originally row 2's data)
$ws.Cells.Item(6, 1).Value = "55"
$ws.Cells.Item(6, 6).Value = 425.7
$ws.Cells.Item(6, 7).Value = 100
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 100
$ws.Cells.Item(6, 10).Value = "08/02/2023"
$ws.Cells.Item(6, 11).Value = "11:26:02.000"
